$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the DSSV (student list) column for all data rows:
# "201612345" -> "20161235"
$ws.Range("E2").Value = "20161234, 20161235, 20161236, 20161237"
$ws.Range("E3").Value = "20161234, 20161235, 20161236, 20161237"
$ws.Range("E4").Value = "20161234, 20161235, 20161236, 20161237"
$ws.Range("E5").Value = "20161234, 20161235, 20161236, 20161237"

# Update the selection to match the author's last selection (cells E4 and E5
# were multi-selected, ending with E5 as the active cell)
$ws.Range("E4").Select()
$ws.Range("E5").Select()
